# Scheduled runner update: refresh Universalis market-board price snapshots
# and recompute leve profit columns across all crafting class sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1178.2
$ws.Range("I2").Value = 1178.2
$ws.Range("K2").Value = 1178.2
$ws.Range("M2").Value = -1065.2

$ws.Range("H15").Value = 2256.6829
$ws.Range("I15").Value = 2256.6829
$ws.Range("K15").Value = 6770.048699999999
$ws.Range("M15").Value = -6601.048699999999

$ws.Range("H33").Value = 239.47826
$ws.Range("I33").Value = 248.09091
$ws.Range("K33").Value = 248.09091
$ws.Range("M33").Value = -19.09091000000001

$ws.Range("H80").Value = 739.5
$ws.Range("I80").Value = 680.75
$ws.Range("J80").Value = 768.875
$ws.Range("K80").Value = 2042.25
$ws.Range("L80").Value = 2306.625
$ws.Range("M80").Value = -1044.25
$ws.Range("N80").Value = -4302.625

$ws.Range("H83").Value = 739.5
$ws.Range("I83").Value = 680.75
$ws.Range("J83").Value = 768.875
$ws.Range("K83").Value = 6126.75
$ws.Range("L83").Value = 6919.875
$ws.Range("M83").Value = -1134.75
$ws.Range("N83").Value = -16903.875

$ws.Range("H98").Value = 760.2222
$ws.Range("I98").Value = 730.25
$ws.Range("K98").Value = 730.25
$ws.Range("M98").Value = 767.75

$ws.Range("H106").Value = 5499.6665
$ws.Range("I106").Value = 4249.5
$ws.Range("J106").Value = 8000
$ws.Range("K106").Value = 4249.5
$ws.Range("L106").Value = 8000
$ws.Range("M106").Value = -3618.5
$ws.Range("N106").Value = -9262

$ws.Range("H122").Value = 760.2222
$ws.Range("I122").Value = 730.25
$ws.Range("K122").Value = 2190.75
$ws.Range("M122").Value = 259.25

$ws.Range("H132").Value = 806.8372000000001
$ws.Range("I132").Value = 802.2381
$ws.Range("K132").Value = 2406.7143
$ws.Range("M132").Value = 123.2856999999999

$ws.Range("H137").Value = 1755.3889
$ws.Range("I137").Value = 1402.48
$ws.Range("K137").Value = 4207.440000000001
$ws.Range("M137").Value = -1657.440000000001

$ws.Range("H138").Value = 4234.8276
$ws.Range("I138").Value = 3599.2
$ws.Range("J138").Value = 4367.25
$ws.Range("K138").Value = 10797.6
$ws.Range("L138").Value = 13101.75
$ws.Range("M138").Value = -5657.599999999999
$ws.Range("N138").Value = -23381.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4803.6196
$ws.Range("I32").Value = 4513.591
$ws.Range("J32").Value = 11184.25
$ws.Range("K32").Value = 4513.591
$ws.Range("L32").Value = 11184.25
$ws.Range("M32").Value = -4226.591
$ws.Range("N32").Value = -11758.25

$ws.Range("H45").Value = 2707
$ws.Range("I45").Value = 2707
$ws.Range("K45").Value = 2707
$ws.Range("M45").Value = -2330

$ws.Range("H63").Value = 750
$ws.Range("I63").Value = 750
$ws.Range("K63").Value = 750
$ws.Range("M63").Value = -64

$ws.Range("H66").Value = 750
$ws.Range("I66").Value = 750
$ws.Range("K66").Value = 3750
$ws.Range("M66").Value = -318

$ws.Range("H74").Value = 1653.4667
$ws.Range("J74").Value = 2402
$ws.Range("L74").Value = 2402
$ws.Range("N74").Value = -4150

$ws.Range("H77").Value = 1653.4667
$ws.Range("J77").Value = 2402
$ws.Range("L77").Value = 12010
$ws.Range("N77").Value = -20746

$ws.Range("H132").Value = 2691.7917
$ws.Range("I132").Value = 1913
$ws.Range("K132").Value = 5739
$ws.Range("M132").Value = -3209

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 4510.5
$ws.Range("I33").Value = 4510.5
$ws.Range("K33").Value = 4510.5
$ws.Range("M33").Value = -4174.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 724.25
$ws.Range("I22").Value = 730.5714
$ws.Range("J22").Value = 680
$ws.Range("K22").Value = 730.5714
$ws.Range("L22").Value = 680
$ws.Range("M22").Value = -380.5714
$ws.Range("N22").Value = -1380

$ws.Range("H35").Value = 500
$ws.Range("I35").Value = 500
$ws.Range("K35").Value = 500
$ws.Range("M35").Value = -206

$ws.Range("H58").Value = 2815.0667
$ws.Range("I58").Value = 2626.8572
$ws.Range("K58").Value = 2626.8572
$ws.Range("M58").Value = -2423.8572

$ws.Range("H86").Value = 22439.723
$ws.Range("J86").Value = 31729.1
$ws.Range("L86").Value = 31729.1
$ws.Range("N86").Value = -33975.1

$ws.Range("H89").Value = 22439.723
$ws.Range("J89").Value = 31729.1
$ws.Range("L89").Value = 158645.5
$ws.Range("N89").Value = -169877.5

$ws.Range("H107").Value = 2971
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 2971
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 2971
$ws.Range("M107").Value = $null
$ws.Range("N107").Value = -6811

$ws.Range("H136").Value = 2815.0667
$ws.Range("I136").Value = 2626.8572
$ws.Range("K136").Value = 7880.571599999999
$ws.Range("M136").Value = -5330.571599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 255249.25
$ws.Range("I80").Value = 4999
$ws.Range("K80").Value = 14997
$ws.Range("M80").Value = -14061

$ws.Range("H83").Value = 255249.25
$ws.Range("I83").Value = 4999
$ws.Range("K83").Value = 44991
$ws.Range("M83").Value = -40311

$ws.Range("H128").Value = 210525.5
$ws.Range("I128").Value = 210525.5
$ws.Range("K128").Value = 631576.5
$ws.Range("M128").Value = -626596.5

$ws.Range("H134").Value = 1632.3334
$ws.Range("I134").Value = 1632.3334
$ws.Range("K134").Value = 4897.0002
$ws.Range("M134").Value = 172.9997999999996

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 428.0357
$ws.Range("I97").Value = 448.14285
$ws.Range("K97").Value = 448.14285
$ws.Range("M97").Value = 47.85714999999999

$ws.Range("H132").Value = 3907.1667
$ws.Range("I132").Value = 3315.1667
$ws.Range("K132").Value = 9945.500100000001
$ws.Range("M132").Value = -7415.500100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").Value = $null

$ws.Range("H132").Value = 5047.25
$ws.Range("I132").Value = 4885.5625
$ws.Range("K132").Value = 14656.6875
$ws.Range("M132").Value = -12126.6875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2953.8667
$ws.Range("I136").Value = 2784.5417
$ws.Range("J136").Value = 3631.1667
$ws.Range("K136").Value = 8353.625100000001
$ws.Range("M136").Value = -5803.625100000001
